# Add a new "Estimates" worksheet after Sheet1, carrying the story-point
# estimate total, and make it the active/visible sheet -- matching the
# commit "Update MOSIP_Partner Management Requirements.xlsx".

$wb = $excel.ActiveWorkbook

# Insert the new sheet right after the existing "Sheet1" so it lands as the
# 2nd sheet (sheetId 2 / rId2), exactly like the target workbook.
$sheet1 = $wb.Worksheets.Item(1)
$ws = $wb.Worksheets.Add($null, $sheet1)
$ws.Name = "Estimates"

# Column widths (B = 56 chars, C ~ 10.5 chars) and content.
$ws.Columns.Item(2).ColumnWidth = 55.1666666
$ws.Columns.Item(3).ColumnWidth = 9.7083333

$ws.Range("B2").Value = "Total Story Points Esimates (Including Desing, Cut Effort, DB Design, Testing, Requirement Detailing, Code Review, Bug Fixing, Documentation, Release Notes)"
$ws.Range("B2").WrapText = $true
$ws.Rows.Item(2).RowHeight = 43.5

$ws.Range("C2").Value = 314

# Page setup / view to match the authored sheet.
$ws.PageSetup.Orientation = 1
$ws.Range("C2").Select() | Out-Null

# Make "Estimates" the active/visible tab (mirrors activeTab="1" + the
# tabSelected move from Sheet1 to Estimates in the diff).
$ws.Activate() | Out-Null
